# Duplicate the language/count data rows (A2:A32 -- i.e. every row except
# the "Languages And Count" header in A1) and append the copy directly
# below the existing data, so the sheet grows from A1:A32 to A1:A63.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A32").Copy()
$ws.Range("A33").PasteSpecial(-4104)
